$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 47657
$ws.Range("J3").Value = 47657
$ws.Range("L3").Value = 47657
$ws.Range("N3").Value = -47885
$ws.Range("H28").Value = 700.6818
$ws.Range("I28").Value = 425.33334
$ws.Range("J28").Value = 1031.1
$ws.Range("K28").Value = 425.33334
$ws.Range("L28").Value = 1031.1
$ws.Range("M28").Value = 59.66665999999998
$ws.Range("N28").Value = -2001.1
$ws.Range("H102").Value = 47657
$ws.Range("J102").Value = 47657
$ws.Range("L102").Value = 47657
$ws.Range("N102").Value = -54147
$ws.Range("H105").Value = 48823.668
$ws.Range("J105").Value = 48823.668
$ws.Range("L105").Value = 48823.668
$ws.Range("N105").Value = -55811.668
$ws.Range("H112").Value = 1472444.2
$ws.Range("I112").Value = 996.6667
$ws.Range("K112").Value = 2990.0001
$ws.Range("M112").Value = -1882.0001
$ws.Range("H126").Value = 46765.332
$ws.Range("J126").Value = 46765.332
$ws.Range("L126").Value = 46765.332
$ws.Range("N126").Value = -56645.332
$ws.Range("H128").Value = 46388
$ws.Range("J128").Value = 46388
$ws.Range("L128").Value = 46388
$ws.Range("N128").Value = -56348
$ws.Range("H133").Value = 54877.555
$ws.Range("J133").Value = 54877.555
$ws.Range("L133").Value = 54877.555
$ws.Range("N133").Value = -64997.555
$ws.Range("H141").Value = 1884.4667
$ws.Range("I141").Value = 935.53845
$ws.Range("J141").Value = 8052.5
$ws.Range("K141").Value = 2806.61535
$ws.Range("L141").Value = 24157.5
$ws.Range("M141").Value = 2373.38465
$ws.Range("N141").Value = -34517.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13376.913
$ws.Range("I32").Value = 13441.216
$ws.Range("K32").Value = 13441.216
$ws.Range("M32").Value = -13154.216
$ws.Range("H80").Value = 38994.668
$ws.Range("J80").Value = 38994.668
$ws.Range("L80").Value = 38994.668
$ws.Range("N80").Value = -40990.668
$ws.Range("H83").Value = 38994.668
$ws.Range("J83").Value = 38994.668
$ws.Range("L83").Value = 116984.004
$ws.Range("N83").Value = -126968.004
$ws.Range("H101").Value = 49598
$ws.Range("J101").Value = 49598
$ws.Range("L101").Value = 49598
$ws.Range("N101").Value = -56088
$ws.Range("H122").Value = 1086.6154
$ws.Range("I122").Value = 956.8889
$ws.Range("J122").Value = 1378.5
$ws.Range("K122").Value = 2870.6667
$ws.Range("L122").Value = 4135.5
$ws.Range("M122").Value = -420.6667000000002
$ws.Range("N122").Value = -9035.5
$ws.Range("H123").Value = 42190.6
$ws.Range("J123").Value = 42190.6
$ws.Range("L123").Value = 42190.6
$ws.Range("N123").Value = -51990.6
$ws.Range("H125").Value = 45331.332
$ws.Range("J125").Value = 45331.332
$ws.Range("L125").Value = 45331.332
$ws.Range("N125").Value = -55171.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 55494.5
$ws.Range("J57").Value = 55494.5
$ws.Range("L57").Value = 55494.5
$ws.Range("N57").Value = -56934.5
$ws.Range("H130").Value = 49889.5
$ws.Range("J130").Value = 49889.5
$ws.Range("L130").Value = 49889.5
$ws.Range("N130").Value = -59929.5
$ws.Range("H133").Value = 40166.5
$ws.Range("J133").Value = 40166.5
$ws.Range("L133").Value = 40166.5
$ws.Range("N133").Value = -50286.5
$ws.Range("H136").Value = 55494.5
$ws.Range("J136").Value = 55494.5
$ws.Range("L136").Value = 55494.5
$ws.Range("N136").Value = -65694.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 73330
$ws.Range("J52").Value = 73330
$ws.Range("L52").Value = 73330
$ws.Range("N52").Value = -73918
$ws.Range("H100").Value = 43436
$ws.Range("J100").Value = 43436
$ws.Range("L100").Value = 43436
$ws.Range("N100").Value = -45600
$ws.Range("H118").Value = 48738
$ws.Range("J118").Value = 48738
$ws.Range("L118").Value = 48738
$ws.Range("N118").Value = -52052
$ws.Range("H137").Value = 45021.355
$ws.Range("J137").Value = 45021.355
$ws.Range("L137").Value = 45021.355
$ws.Range("N137").Value = -55221.355
$ws.Range("H139").Value = 61299.8
$ws.Range("J139").Value = 66124.75
$ws.Range("L139").Value = 66124.75
$ws.Range("N139").Value = -76404.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 3280.7058
$ws.Range("J115").Value = 3584.8
$ws.Range("L115").Value = 10754.4
$ws.Range("N115").Value = -13104.4
$ws.Range("I134").Value = 47669310
$ws.Range("K134").Value = 143007930
$ws.Range("M134").Value = -143002860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5109.524
$ws.Range("I70").Value = 5133.3335
$ws.Range("J70").Value = 4966.6665
$ws.Range("K70").Value = 5133.3335
$ws.Range("L70").Value = 4966.6665
$ws.Range("M70").Value = -4863.3335
$ws.Range("N70").Value = -5506.6665
$ws.Range("H73").Value = 5109.524
$ws.Range("I73").Value = 5133.3335
$ws.Range("J73").Value = 4966.6665
$ws.Range("K73").Value = 5133.3335
$ws.Range("L73").Value = 4966.6665
$ws.Range("M73").Value = -4197.3335
$ws.Range("N73").Value = -6838.6665
$ws.Range("H122").Value = 1200
$ws.Range("I122").Value = 1225
$ws.Range("K122").Value = 3675
$ws.Range("M122").Value = -1225
$ws.Range("H124").Value = 37920.668
$ws.Range("J124").Value = 37920.668
$ws.Range("L124").Value = 37920.668
$ws.Range("N124").Value = -47740.668
$ws.Range("H126").Value = 4230.6
$ws.Range("I126").Value = 4317.6665
$ws.Range("J126").Value = 4100
$ws.Range("K126").Value = 12952.9995
$ws.Range("L126").Value = 12300
$ws.Range("M126").Value = -10482.9995
$ws.Range("N126").Value = -17240
$ws.Range("H132").Value = 3282.8635
$ws.Range("I132").Value = 1761.7
$ws.Range("J132").Value = 4550.5
$ws.Range("K132").Value = 5285.1
$ws.Range("L132").Value = 13651.5
$ws.Range("M132").Value = -2755.1
$ws.Range("N132").Value = -18711.5
$ws.Range("H135").Value = 49675
$ws.Range("J135").Value = 49675
$ws.Range("L135").Value = 49675
$ws.Range("N135").Value = -59815
$ws.Range("H138").Value = 53500
$ws.Range("J138").Value = 53500
$ws.Range("L138").Value = 53500
$ws.Range("N138").Value = -63780

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 48670
$ws.Range("J36").Value = 48670
$ws.Range("L36").Value = 48670
$ws.Range("N36").Value = -49794
$ws.Range("H40").Value = 4746.5
$ws.Range("I40").Value = 3000.3333
$ws.Range("J40").Value = 9985
$ws.Range("K40").Value = 3000.3333
$ws.Range("L40").Value = 9985
$ws.Range("M40").Value = -2864.3333
$ws.Range("N40").Value = -10257
$ws.Range("H124").Value = 48429
$ws.Range("J124").Value = 48429
$ws.Range("L124").Value = 48429
$ws.Range("N124").Value = -58249
$ws.Range("H134").Value = 50728.25
$ws.Range("J134").Value = 50728.25
$ws.Range("L134").Value = 50728.25
$ws.Range("N134").Value = -60868.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 53795
$ws.Range("J46").Value = 53795
$ws.Range("L46").Value = 53795
$ws.Range("N46").Value = -54257
$ws.Range("H122").Value = 4762680.5
$ws.Range("I122").Value = 5715016.5
$ws.Range("K122").Value = 17145049.5
$ws.Range("M122").Value = -17142599.5
$ws.Range("H131").Value = 49116.832
$ws.Range("J131").Value = 49116.832
$ws.Range("L131").Value = 49116.832
$ws.Range("N131").Value = -59196.832
$ws.Range("H134").Value = 53795
$ws.Range("J134").Value = 53795
$ws.Range("L134").Value = 161385
$ws.Range("N134").Value = -166455

Write-Host "Applied all updates"
